$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A34").Value = "Gennaro Bullo"
$ws.Range("B34").Value = "ELIA BATTISTI | U.S. Guarna"
$ws.Range("C34").Value = "ENRICO BORDIGNON | Pinguini Trentini"
$ws.Range("D34").Value = "Leonardo Viola | Shark Attack"
$ws.Range("E34").Value = "Randy Cobbinah | MAI UNA GIOIA"
$ws.Range("F34").Value = "Matteo Mazzola | MediaserT"
